# Updated symbol list on Mon Dec 19 10:48:21 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

# Column D (Price) updates - values kept as text to preserve exact formatting
Set-TextValue "D2" "248.44"
Set-TextValue "D3" "21.71"
Set-TextValue "D4" "5.481"
Set-TextValue "D5" "0.05697"
Set-TextValue "D6" "3.367"
Set-TextValue "D7" "0.8038"
Set-TextValue "D9" "0.1507"
Set-TextValue "D10" "0.07373"
Set-TextValue "D12" "0.03014"
Set-TextValue "D13" "0.09287"
Set-TextValue "D14" "3.434"
Set-TextValue "D15" "0.001645"
Set-TextValue "D16" "0.04695"
Set-TextValue "D17" "0.0005858"

$ws.Range("E17").Value = "16OneONEWorstin24h"

Set-TextValue "D19" "0.005056"
Set-TextValue "D20" "0.001043"
Set-TextValue "D21" "0.0001500"
Set-TextValue "D22" "0.0003199"
Set-TextValue "D23" "3.767"
Set-TextValue "D24" "6.427"
Set-TextValue "D25" "2.130"

Set-TextValue "D40" "0.04119"
Set-TextValue "D41" "0.006945"

# Rows 42 and 43 swap coin data (BKEXToken <-> CEJI)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003500"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1044"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue "D44" "0.009142"
Set-TextValue "D45" "0.00005839"

Set-TextValue "D47" "0.0005498"
$ws.Range("E47").Value = "46ACDXExchangeACXT"

Set-TextValue "D48" "0.6822"
Set-TextValue "D49" "0.009266"
